# InitProperty.xlsx — "unify the conception of DataNode, DataTable, Entity"
#
# The only content-level edits in this commit are:
#   1) The sheet formerly named "Property1" is renamed to "DataNode"
#      (part of a repo-wide rename unifying DataNode/DataTable/Entity
#      naming across the Excel config files).
#   2) The author's view/selection moved: the frozen pane was scrolled
#      further down and the active cell ended on E50.
#
# Everything else in the raw XML diff (fileVersion/rupBuild bump, xr2/xr/
# xr3/x16r2/mc:Ignorable namespace churn, the absPath of the machine that
# re-saved the file, font substitution Calibri -> 宋体, default row
# height/column width recalculation, phoneticPr, extra xr:uid attributes,
# timeline style extension, etc.) is mechanical fallout of the workbook
# being re-saved by a different Excel build/locale and isn't an
# intentional edit, so it is not something to reproduce cell-by-cell here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet: Property1 -> DataNode
$ws.Name = "DataNode"

# 2) Move the selection to where the author left it (E50). The engine
#    re-derives activeCell/sqref for the sheetView's selection from the
#    current selected range.
$ws.Range("E50").Select()
